$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.793.81"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.104.32"
$ws.Range("E3").Value = "  +0.62%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.50"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("E6").Value = "  +1.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.86"
$ws.Range("E7").Value = "  +2.42%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +1.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0843"
$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("E11").Value = "  -0.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.86"
$ws.Range("E12").Value = "  +6.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.414.70"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.97"
$ws.Range("E14").Value = "  -0.36%  "

$ws.Range("E15").Value = "  +0.56%  "

$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.095.78"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.820.97"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.62"
$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.03"
$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("E21").Value = "  +1.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.96"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("E24").Value = "  -0.88%  "

$ws.Range("E25").Value = "  -0.27%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.65"
$ws.Range("E26").Value = "  +1.64%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.34"
$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("E29").Value = "  +0.78%  "

$ws.Range("E30").Value = "  +1.07%  "

$ws.Range("E32").Value = "  +0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.57"
$ws.Range("E33").Value = "  +2.18%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.78"
$ws.Range("E34").Value = "  +1.12%  "

$ws.Range("B35").Value = "THORChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.14"
$ws.Range("E35").Value = "  +11.64%  "

$ws.Range("E36").Value = "  +0.36%  "

$ws.Range("E37").Value = "  -0.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.50"
$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("E40").Value = "  +3.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.94"
$ws.Range("E41").Value = "  -2.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.00"
$ws.Range("E42").Value = "  +1.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.528.64"
$ws.Range("E43").Value = "  -0.58%  "

$ws.Range("E44").Value = "  +7.56%  "

$ws.Range("E45").Value = "  -0.74%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0914"
$ws.Range("E46").Value = "  -1.31%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.77"
$ws.Range("E47").Value = "  +1.68%  "

$ws.Range("E48").Value = "  +5.23%  "

$ws.Range("E49").Value = "  +0.93%  "

$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.301.45"
$ws.Range("E51").Value = "  +0.55%  "
